$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.50081833333333
$ws.Range("H2").Value = 46.502455
$ws.Range("I2").Value = 0.6975713722222094
$ws.Range("J2").Value = 0.6975713722222093
$ws.Range("M2").Value = 29.01761566666667
$ws.Range("N2").Value = 87.052847
$ws.Range("O2").Value = 0.6436075952942075
$ws.Range("P2").Value = 0.6436075952942075
$ws.Range("Q2").Value = 449.7967889154872
$ws.Range("R2").Value = 4048.171100239385
$ws.Range("S2").Value = 0.4489622334220167
$ws.Range("T2").Value = 0.4489622334220167
$ws.Range("G3").Value = 15.50081833333333
$ws.Range("H3").Value = 46.502455
$ws.Range("I3").Value = 0.6975713722222094
$ws.Range("J3").Value = 0.6975713722222093
$ws.Range("O3").Value = 0.05796185537580412
$ws.Range("P3").Value = 0.05796185537580412
$ws.Range("Q3").Value = 40.50768918552445
$ws.Range("R3").Value = 364.56920266972
$ws.Range("S3").Value = 0.04043253099104493
$ws.Range("T3").Value = 0.04043253099104492
$ws.Range("G4").Value = 15.50081833333333
$ws.Range("H4").Value = 46.502455
$ws.Range("I4").Value = 0.6975713722222094
$ws.Range("J4").Value = 0.6975713722222093
$ws.Range("M4").Value = 1.123006
$ws.Range("N4").Value = 3.369018
$ws.Range("O4").Value = 0.02490815232594174
$ws.Range("P4").Value = 0.02490815232594174
$ws.Range("Q4").Value = 17.40751199324333
$ws.Range("R4").Value = 156.66760793919
$ws.Range("S4").Value = 0.017375213997527
$ws.Range("T4").Value = 0.01737521399752699
$ws.Range("G5").Value = 15.50081833333333
$ws.Range("H5").Value = 46.502455
$ws.Range("I5").Value = 0.6975713722222094
$ws.Range("J5").Value = 0.6975713722222093
$ws.Range("M5").Value = 10.33196133333333
$ws.Range("N5").Value = 30.995884
$ws.Range("O5").Value = 0.2291617914030796
$ws.Range("P5").Value = 0.2291617914030796
$ws.Range("Q5").Value = 160.1538556550244
$ws.Range("R5").Value = 1441.38470089522
$ws.Range("S5").Value = 0.1598567052899459
$ws.Range("T5").Value = 0.1598567052899459
$ws.Range("G6").Value = 15.50081833333333
$ws.Range("H6").Value = 46.502455
$ws.Range("I6").Value = 0.6975713722222094
$ws.Range("J6").Value = 0.6975713722222093
$ws.Range("M6").Value = 2.000037
$ws.Range("N6").Value = 6.000111
$ws.Range("O6").Value = 0.044360605600967
$ws.Range("P6").Value = 0.044360605600967
$ws.Range("Q6").Value = 31.002210196945
$ws.Range("R6").Value = 279.019891772505
$ws.Range("S6").Value = 0.03094468852167478
$ws.Range("T6").Value = 0.03094468852167478
$ws.Range("I7").Value = 0.1333925762673397
$ws.Range("J7").Value = 0.1333925762673396
$ws.Range("M7").Value = 29.01761566666667
$ws.Range("N7").Value = 87.052847
$ws.Range("O7").Value = 0.6436075952942075
$ws.Range("P7").Value = 0.6436075952942075
$ws.Range("Q7").Value = 86.01206250634512
$ws.Range("R7").Value = 774.108562557106
$ws.Range("S7").Value = 0.08585247524152166
$ws.Range("T7").Value = 0.08585247524152163
$ws.Range("I8").Value = 0.1333925762673397
$ws.Range("J8").Value = 0.1333925762673396
$ws.Range("O8").Value = 0.05796185537580412
$ws.Range("P8").Value = 0.05796185537580412
$ws.Range("S8").Value = 0.007731681213813463
$ws.Range("T8").Value = 0.007731681213813462
$ws.Range("I9").Value = 0.1333925762673397
$ws.Range("J9").Value = 0.1333925762673396
$ws.Range("M9").Value = 1.123006
$ws.Range("N9").Value = 3.369018
$ws.Range("O9").Value = 0.02490815232594174
$ws.Range("P9").Value = 0.02490815232594174
$ws.Range("Q9").Value = 3.328738769462666
$ws.Range("R9").Value = 29.958648925164
$ws.Range("S9").Value = 0.003322562608816697
$ws.Range("T9").Value = 0.003322562608816697
$ws.Range("I10").Value = 0.1333925762673397
$ws.Range("J10").Value = 0.1333925762673396
$ws.Range("M10").Value = 10.33196133333333
$ws.Range("N10").Value = 30.995884
$ws.Range("O10").Value = 0.2291617914030796
$ws.Range("P10").Value = 0.2291617914030796
$ws.Range("Q10").Value = 30.62530409887022
$ws.Range("R10").Value = 275.627736889832
$ws.Range("S10").Value = 0.03056848173729548
$ws.Range("T10").Value = 0.03056848173729547
$ws.Range("I11").Value = 0.1333925762673397
$ws.Range("J11").Value = 0.1333925762673396
$ws.Range("M11").Value = 2.000037
$ws.Range("N11").Value = 6.000111
$ws.Range("O11").Value = 0.044360605600967
$ws.Range("P11").Value = 0.044360605600967
$ws.Range("Q11").Value = 5.928375006242001
$ws.Range("R11").Value = 53.355375056178
$ws.Range("S11").Value = 0.005917375465892366
$ws.Range("T11").Value = 0.005917375465892365
$ws.Range("G12").Value = 1.908787666666667
$ws.Range("H12").Value = 5.726363
$ws.Range("I12").Value = 0.08589969918260204
$ws.Range("J12").Value = 0.08589969918260203
$ws.Range("M12").Value = 29.01761566666667
$ws.Range("N12").Value = 87.052847
$ws.Range("O12").Value = 0.6436075952942075
$ws.Range("P12").Value = 0.6436075952942075
$ws.Range("Q12").Value = 55.38846690060678
$ws.Range("R12").Value = 498.496202105461
$ws.Range("S12").Value = 0.0552856988274103
$ws.Range("T12").Value = 0.05528569882741029
$ws.Range("G13").Value = 1.908787666666667
$ws.Range("H13").Value = 5.726363
$ws.Range("I13").Value = 0.08589969918260204
$ws.Range("J13").Value = 0.08589969918260203
$ws.Range("O13").Value = 0.05796185537580412
$ws.Range("P13").Value = 0.05796185537580412
$ws.Range("Q13").Value = 4.988161002843556
$ws.Range("R13").Value = 44.893449025592
$ws.Range("S13").Value = 0.00497890594084706
$ws.Range("T13").Value = 0.004978905940847059
$ws.Range("G14").Value = 1.908787666666667
$ws.Range("H14").Value = 5.726363
$ws.Range("I14").Value = 0.08589969918260204
$ws.Range("J14").Value = 0.08589969918260203
$ws.Range("M14").Value = 1.123006
$ws.Range("N14").Value = 3.369018
$ws.Range("O14").Value = 0.02490815232594174
$ws.Range("P14").Value = 0.02490815232594174
$ws.Range("Q14").Value = 2.143580002392667
$ws.Range("R14").Value = 19.292220021534
$ws.Range("S14").Value = 0.002139602791992825
$ws.Range("T14").Value = 0.002139602791992824
$ws.Range("G15").Value = 1.908787666666667
$ws.Range("H15").Value = 5.726363
$ws.Range("I15").Value = 0.08589969918260204
$ws.Range("J15").Value = 0.08589969918260203
$ws.Range("M15").Value = 10.33196133333333
$ws.Range("N15").Value = 30.995884
$ws.Range("O15").Value = 0.2291617914030796
$ws.Range("P15").Value = 0.2291617914030796
$ws.Range("Q15").Value = 19.72152036554355
$ws.Range("R15").Value = 177.493683289892
$ws.Range("S15").Value = 0.01968492894567073
$ws.Range("T15").Value = 0.01968492894567073
$ws.Range("G16").Value = 1.908787666666667
$ws.Range("H16").Value = 5.726363
$ws.Range("I16").Value = 0.08589969918260204
$ws.Range("J16").Value = 0.08589969918260203
$ws.Range("M16").Value = 2.000037
$ws.Range("N16").Value = 6.000111
$ws.Range("O16").Value = 0.044360605600967
$ws.Range("P16").Value = 0.044360605600967
$ws.Range("Q16").Value = 3.817645958477001
$ws.Range("R16").Value = 34.35881362629301
$ws.Range("S16").Value = 0.003810562676681117
$ws.Range("T16").Value = 0.003810562676681116
$ws.Range("G17").Value = 0.5887749999999999
$ws.Range("H17").Value = 1.766325
$ws.Range("I17").Value = 0.02649618722367226
$ws.Range("J17").Value = 0.02649618722367225
$ws.Range("M17").Value = 29.01761566666667
$ws.Range("N17").Value = 87.052847
$ws.Range("O17").Value = 0.6436075952942075
$ws.Range("P17").Value = 0.6436075952942075
$ws.Range("Q17").Value = 17.08484666414166
$ws.Range("R17").Value = 153.763619977275
$ws.Range("S17").Value = 0.0170531473434928
$ws.Range("T17").Value = 0.0170531473434928
$ws.Range("G18").Value = 0.5887749999999999
$ws.Range("H18").Value = 1.766325
$ws.Range("I18").Value = 0.02649618722367226
$ws.Range("J18").Value = 0.02649618722367225
$ws.Range("O18").Value = 0.05796185537580412
$ws.Range("P18").Value = 0.05796185537580412
$ws.Range("Q18").Value = 1.538622941533333
$ws.Range("R18").Value = 13.8476064738
$ws.Range("S18").Value = 0.00153576817186872
$ws.Range("T18").Value = 0.00153576817186872
$ws.Range("G19").Value = 0.5887749999999999
$ws.Range("H19").Value = 1.766325
$ws.Range("I19").Value = 0.02649618722367226
$ws.Range("J19").Value = 0.02649618722367225
$ws.Range("M19").Value = 1.123006
$ws.Range("N19").Value = 3.369018
$ws.Range("O19").Value = 0.02490815232594174
$ws.Range("P19").Value = 0.02490815232594174
$ws.Range("Q19").Value = 0.6611978576499999
$ws.Range("R19").Value = 5.950780718849999
$ws.Range("S19").Value = 0.0006599710674239
$ws.Range("T19").Value = 0.0006599710674238998
$ws.Range("G20").Value = 0.5887749999999999
$ws.Range("H20").Value = 1.766325
$ws.Range("I20").Value = 0.02649618722367226
$ws.Range("J20").Value = 0.02649618722367225
$ws.Range("M20").Value = 10.33196133333333
$ws.Range("N20").Value = 30.995884
$ws.Range("O20").Value = 0.2291617914030796
$ws.Range("P20").Value = 0.2291617914030796
$ws.Range("Q20").Value = 6.083200534033332
$ws.Range("R20").Value = 54.74880480629999
$ws.Range("S20").Value = 0.006071913729528124
$ws.Range("T20").Value = 0.006071913729528123
$ws.Range("G21").Value = 0.5887749999999999
$ws.Range("H21").Value = 1.766325
$ws.Range("I21").Value = 0.02649618722367226
$ws.Range("J21").Value = 0.02649618722367225
$ws.Range("M21").Value = 2.000037
$ws.Range("N21").Value = 6.000111
$ws.Range("O21").Value = 0.044360605600967
$ws.Range("P21").Value = 0.044360605600967
$ws.Range("Q21").Value = 1.177571784675
$ws.Range("R21").Value = 10.598146062075
$ws.Range("S21").Value = 0.001175386911358706
$ws.Range("T21").Value = 0.001175386911358706
$ws.Range("G22").Value = 1.258608
$ws.Range("H22").Value = 3.775824
$ws.Range("I22").Value = 0.0566401651041768
$ws.Range("J22").Value = 0.05664016510417678
$ws.Range("M22").Value = 29.01761566666667
$ws.Range("N22").Value = 87.052847
$ws.Range("O22").Value = 0.6436075952942075
$ws.Range("P22").Value = 0.6436075952942075
$ws.Range("Q22").Value = 36.521803218992
$ws.Range("R22").Value = 328.696228970928
$ws.Range("S22").Value = 0.03645404045976611
$ws.Range("T22").Value = 0.03645404045976611
$ws.Range("G23").Value = 1.258608
$ws.Range("H23").Value = 3.775824
$ws.Range("I23").Value = 0.0566401651041768
$ws.Range("J23").Value = 0.05664016510417678
$ws.Range("O23").Value = 0.05796185537580412
$ws.Range("P23").Value = 0.05796185537580412
$ws.Range("Q23").Value = 3.289071620224
$ws.Range("R23").Value = 29.601644582016
$ws.Range("S23").Value = 0.003282969058229963
$ws.Range("T23").Value = 0.003282969058229962
$ws.Range("G24").Value = 1.258608
$ws.Range("H24").Value = 3.775824
$ws.Range("I24").Value = 0.0566401651041768
$ws.Range("J24").Value = 0.05664016510417678
$ws.Range("M24").Value = 1.123006
$ws.Range("N24").Value = 3.369018
$ws.Range("O24").Value = 0.02490815232594174
$ws.Range("P24").Value = 0.02490815232594174
$ws.Range("Q24").Value = 1.413424335648
$ws.Range("R24").Value = 12.720819020832
$ws.Range("S24").Value = 0.001410801860181326
$ws.Range("T24").Value = 0.001410801860181325
$ws.Range("G25").Value = 1.258608
$ws.Range("H25").Value = 3.775824
$ws.Range("I25").Value = 0.0566401651041768
$ws.Range("J25").Value = 0.05664016510417678
$ws.Range("M25").Value = 10.33196133333333
$ws.Range("N25").Value = 30.995884
$ws.Range("O25").Value = 0.2291617914030796
$ws.Range("P25").Value = 0.2291617914030796
$ws.Range("Q25").Value = 13.003889189824
$ws.Range("R25").Value = 117.035002708416
$ws.Range("S25").Value = 0.01297976170063935
$ws.Range("T25").Value = 0.01297976170063935
$ws.Range("G26").Value = 1.258608
$ws.Range("H26").Value = 3.775824
$ws.Range("I26").Value = 0.0566401651041768
$ws.Range("J26").Value = 0.05664016510417678
$ws.Range("M26").Value = 2.000037
$ws.Range("N26").Value = 6.000111
$ws.Range("O26").Value = 0.044360605600967
$ws.Range("P26").Value = 0.044360605600967
$ws.Range("Q26").Value = 2.517262568496
$ws.Range("R26").Value = 22.655363116464
$ws.Range("S26").Value = 0.002512592025360041
$ws.Range("T26").Value = 0.002512592025360041
